$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shorten the header labels (drop units, now captured/implied elsewhere)
$ws.Range("B1").Value = "wt"
$ws.Range("C1").Value = "Onset Temp"
$ws.Range("D1").Value = "Weight Loss"
$ws.Range("E1").Value = "DTG"
$ws.Range("F1").Value = "Non-Volatile Residue"
$ws.Range("G1").Value = "Glass Transition Temp"

# Update selected cell
$ws.Range("B2").Select()
